$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2:G15").ClearContents() | Out-Null
$ws.Range("H14").Select() | Out-Null
